$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged) - update B3, C3, D3
$ws.Range("B3").Value = 0.9935514042246122
$ws.Range("C3").Value = 0.9934281601580318
$ws.Range("D3").Value = 0.9920871004023777

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update B4, C4, D4
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9934653510987963
$ws.Range("C4").Value = 0.9939682548970614
$ws.Range("D4").Value = 0.9939864265416783

# Row 5: AdaBoostRegressor -> MLPRegressor, update B5, C5, D5
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9949739930920797
$ws.Range("C5").Value = 0.9954749373589905
$ws.Range("D5").Value = 0.9959084756909334
